# Generate Report for Archive
#
# 1) Every "Ready for handoff" status value becomes "In Translation"
#    (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 all share that string).
# 2) The Status-related columns got narrower (Overview columns E & F,
#    and column C on the two language sheets).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the literal string on the LEFT of -eq. PowerShell picks
        # the comparison type from the left-hand operand; with a boolean
        # cell (Value2 = $true) on the left, "-eq $oldStatus" would coerce
        # the string to a bool (any non-empty string -> $true) and produce
        # false positives (e.g. the "To be localized" True/False column).
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value2 = $newStatus
        }
    }
}

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
